# Se arregla la manera en que se guardan las agendas disponibles desde el
# perfil de docente y mentor. Este script reproduce los nuevos registros
# agregados al libro: una franja horaria disponible en "Agendas", tres
# nuevos estudiantes registrados y una nueva asesoria agendada.
#
# Nota: algunos valores (documentos, telefonos, contraseñas numericas y la
# fecha "02-12-2023") lucen como numeros, por lo que se fuerza su formato
# a texto ("@") antes de escribirlos para que se guarden como cadenas,
# igual que el resto de los datos de la hoja, y luego se restaura el
# estilo normal de la celda.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# 1) Nueva franja horaria disponible (hoja "Agendas")
$wsAgendas = $wb.Worksheets.Item("Agendas")
$wsAgendas.Range("A2").Value = "s"
$wsAgendas.Range("B2").Value = "09:00-13:00"

# 2) Nuevo estudiante: juandiego mesa (hoja "estudiantes", fila 4)
$wsEst = $wb.Worksheets.Item("estudiantes")
$wsEst.Range("A4").Value = "juandiego"
$wsEst.Range("B4").Value = "mesa"
Set-TextValue $wsEst.Range("C4") "1001469998"
$wsEst.Range("D4").Value = "hades666"
$wsEst.Range("E4").Value = "juandiegomesa234@gmaul.com"
Set-TextValue $wsEst.Range("F4") "3226505292"
Set-TextValue $wsEst.Range("G4") "123456"
$wsEst.Range("H4").Value = "Estudiante"
$wsEst.Range("I4").Value = "Antioquia"
$wsEst.Range("J4").Value = "Medellín"
$wsEst.Range("K4").Value = "P.C.J.I.C"

# 3) Nuevo estudiante: Mel Suarez (hoja "estudiantes", fila 5)
$wsEst.Range("A5").Value = "Mel"
$wsEst.Range("B5").Value = "Suarez"
Set-TextValue $wsEst.Range("C5") "1011392080"
$wsEst.Range("D5").Value = "Mel1227"
$wsEst.Range("E5").Value = "mel@gmail.com"
Set-TextValue $wsEst.Range("F5") "3205727115"
Set-TextValue $wsEst.Range("G5") "1234"
$wsEst.Range("H5").Value = "Estudiante"
$wsEst.Range("I5").Value = "Antioquia"
$wsEst.Range("J5").Value = "Medellín"
$wsEst.Range("K5").Value = "P.C.J.I.C"

# 4) Nueva asesoria agendada por Mel Suarez con Maryem Ruiz (hoja "asesorias", fila 3)
$wsAses = $wb.Worksheets.Item("asesorias")
$wsAses.Range("A3").Value = "Mel Suarez"
$wsAses.Range("B3").Value = "Mel1227"
$wsAses.Range("C3").Value = "s"
$wsAses.Range("D3").Value = "Maryem Ruiz"
$wsAses.Range("E3").Value = "Asesoría PPI"
Set-TextValue $wsAses.Range("F3") "02-12-2023"
$wsAses.Range("G3").Value = "06:20 - 06:40"

# 5) Nuevo estudiante: John Garcia (hoja "estudiantes", fila 6)
$wsEst.Range("A6").Value = "John"
$wsEst.Range("B6").Value = "Garcia"
Set-TextValue $wsEst.Range("C6") "1000534371"
$wsEst.Range("D6").Value = "john78tigre"
$wsEst.Range("E6").Value = "john78tigres@gmail.com"
Set-TextValue $wsEst.Range("F6") "3108929832"
Set-TextValue $wsEst.Range("G6") "12345678"
$wsEst.Range("H6").Value = "Estudiante"
$wsEst.Range("I6").Value = "Antioquia"
$wsEst.Range("J6").Value = "Medellín"
$wsEst.Range("K6").Value = "P.C.J.I.C"

$wb.Save()
